$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-detected as a number by Excel, so they stay text like the rest
# of the (inline-string) column.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.668.43'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.597.52'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '211.59'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').Value = '0.513'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = '1.821.64'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '1.614.66'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').Value = '26.644.66'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '0.0₃0751'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('D19').Value = '210.28'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E21').Value = '  +4.38%  '
$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').Value = '8.97'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').Value = '15.35'
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('D34').Value = '1.287.10'
$ws.Range('E34').Value = '  -0.72%  '
$ws.Range('D35').Value = '0.619'
$ws.Range('E35').Value = '  -6.30%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '1.06'
$ws.Range('E39').Value = '  +17.86%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '0.827'
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('D41').Value = '5.45'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').Value = '0.782'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('D44').Value = '63.31'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').Value = '1.734.25'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').Value = '91.34'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('D47').Value = '1.56'
$ws.Range('E47').Value = '  -2.81%  '
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = '7.32'
$ws.Range('E51').Value = '  -2.59%  '

# Drop back to the default (unstyled) cell style now that the text is
# committed, so no stray style index lingers on these cells.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
